$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.933.22'
$ws.Range("E2").Value = '  +2.95%  '
$ws.Range("D3").Value = '1.911.94'
$ws.Range("E3").Value = '  +1.47%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5015'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06874'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.85%  '
$ws.Range("D10").Value = '1.908.53'
$ws.Range("E10").Value = '  +1.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.05'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07326'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '91.46'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.115'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6841'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.96%  '
$ws.Range("D16").Value = '30.910.61'
$ws.Range("E16").Value = '  +2.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008062'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '2.153.66'
$ws.Range("E20").Value = '  +1.56%  '
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("E22").Value = '  +2.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '183.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +34.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.130'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.393'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.952'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.41%  '
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.360'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08997'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.068'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05277'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7500'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.63%  '
$ws.Range("E35").Value = '  +3.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.669'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01955'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +18.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.742'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.187'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9377'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4425'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.871'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.46%  '
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.796'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1348'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05850'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.3943'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.627'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.26%  '
$ws.Range("E51").Value = '  +4.40%  '
